# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.334.32"
$ws.Range("E2").Value = "  -5.66%  "
$ws.Range("D3").Value = "2.225.71"
$ws.Range("E3").Value = "  -5.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.03%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.11%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.23%  "
$ws.Range("D15").Value = "2.556.71"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.28%  "
$ws.Range("E17").Value = "  -8.55%  "
$ws.Range("D18").Value = "2.220.74"
$ws.Range("E18").Value = "  -5.59%  "
$ws.Range("D19").Value = "41.314.72"
$ws.Range("E19").Value = "  -5.51%  "
$ws.Range("D20").Value = "0.0₃0956"
$ws.Range("E20").Value = "  -7.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.77%  "
$ws.Range("E22").Value = "  -7.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.34%  "
$ws.Range("E24").Value = "  +11.88%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.40%  "
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.70%  "
$ws.Range("E32").Value = "  -7.77%  "
$ws.Range("E33").Value = "  -6.94%  "
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("E36").Value = "  -10.24%  "
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0276"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("E50").Value = "  -6.08%  "
$ws.Range("E51").Value = "  -4.61%  "
